$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 132
$ws.Range("H132").Value = 913.36664
$ws.Range("I132").Value = 589.94446
$ws.Range("J132").Value = 3824.1667
$ws.Range("K132").Value = 1769.83338
$ws.Range("L132").Value = 11472.5001
$ws.Range("M132").Value = 760.16662
$ws.Range("N132").Value = -16532.5001

# Row 133
$ws.Range("H133").Value = 34885
$ws.Range("J133").Value = 34885
$ws.Range("L133").Value = 34885
$ws.Range("N133").Value = -45005

# Row 134
$ws.Range("H134").Value = 41890
$ws.Range("J134").Value = 41890
$ws.Range("L134").Value = 41890
$ws.Range("N134").Value = -52030

# Row 136
$ws.Range("H136").Value = 42570
$ws.Range("J136").Value = 42570
$ws.Range("L136").Value = 42570
$ws.Range("N136").Value = -52770

# Row 137
$ws.Range("H137").Value = 1551.9375
$ws.Range("I137").Value = 1327.5834
$ws.Range("J137").Value = 2225
$ws.Range("K137").Value = 3982.7502
$ws.Range("L137").Value = 6675
$ws.Range("M137").Value = -1432.7502
$ws.Range("N137").Value = -11775

# Row 139
$ws.Range("H139").Value = 45773.332
$ws.Range("J139").Value = 45773.332
$ws.Range("L139").Value = 45773.332
$ws.Range("N139").Value = -56053.332

# Row 140
$ws.Range("H140").Value = 49765
$ws.Range("J140").Value = 49765
$ws.Range("L140").Value = 49765
$ws.Range("N140").Value = -60125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 7004.907
$ws.Range("I32").Value = 5889.2896
$ws.Range("J32").Value = 11042.381
$ws.Range("K32").Value = 5889.2896
$ws.Range("L32").Value = 11042.381
$ws.Range("M32").Value = -5602.2896
$ws.Range("N32").Value = -11616.381

# Row 61
$ws.Range("H61").Value = 3908.7317
$ws.Range("I61").Value = 3914.7222
$ws.Range("J61").Value = 3865.6
$ws.Range("K61").Value = 3914.7222
$ws.Range("L61").Value = 3865.6
$ws.Range("M61").Value = -3702.7222
$ws.Range("N61").Value = -4289.6

# Row 74
$ws.Range("H74").Value = 1100.6323
$ws.Range("I74").Value = 840.7193
$ws.Range("J74").Value = 2447.4546
$ws.Range("K74").Value = 840.7193
$ws.Range("L74").Value = 2447.4546
$ws.Range("M74").Value = 33.28070000000002
$ws.Range("N74").Value = -4195.4546

# Row 77
$ws.Range("H77").Value = 1100.6323
$ws.Range("I77").Value = 840.7193
$ws.Range("J77").Value = 2447.4546
$ws.Range("K77").Value = 4203.5965
$ws.Range("L77").Value = 12237.273
$ws.Range("M77").Value = 164.4035000000003
$ws.Range("N77").Value = -20973.273

# Row 132
$ws.Range("H132").Value = 2776.6736
$ws.Range("I132").Value = 1799.2894
$ws.Range("J132").Value = 6153.091
$ws.Range("K132").Value = 5397.8682
$ws.Range("L132").Value = 18459.273
$ws.Range("M132").Value = -2867.8682
$ws.Range("N132").Value = -23519.273

# Row 135
$ws.Range("H135").Value = 49943
$ws.Range("J135").Value = 49943
$ws.Range("L135").Value = 49943
$ws.Range("N135").Value = -60083

# Row 136
$ws.Range("H136").Value = 3908.7317
$ws.Range("I136").Value = 3914.7222
$ws.Range("J136").Value = 3865.6
$ws.Range("K136").Value = 11744.1666
$ws.Range("L136").Value = 11596.8
$ws.Range("M136").Value = -9194.1666
$ws.Range("N136").Value = -16696.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 81
$ws.Range("H81").Value = 32350
$ws.Range("J81").Value = 32350
$ws.Range("L81").Value = 32350
$ws.Range("N81").Value = -34472

# Row 84
$ws.Range("H84").Value = 32350
$ws.Range("J84").Value = 32350
$ws.Range("L84").Value = 97050
$ws.Range("N84").Value = -107658

# Row 134
$ws.Range("H134").Value = 3202.7534
$ws.Range("I134").Value = 3455.2666
$ws.Range("J134").Value = 2796.9285
$ws.Range("K134").Value = 10365.7998
$ws.Range("L134").Value = 8390.7855
$ws.Range("M134").Value = -7830.799800000001
$ws.Range("N134").Value = -13460.7855

# Row 135
$ws.Range("H135").Value = 39933.77
$ws.Range("J135").Value = 39933.77
$ws.Range("L135").Value = 39933.77
$ws.Range("N135").Value = -50073.77

# Row 137
$ws.Range("H137").Value = 62333.332
$ws.Range("J137").Value = 38500
$ws.Range("L137").Value = 38500
$ws.Range("N137").Value = -48700

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 58
$ws.Range("H58").Value = 1399.78
$ws.Range("I58").Value = 766.4838999999999
$ws.Range("J58").Value = 2433.0527
$ws.Range("K58").Value = 766.4838999999999
$ws.Range("L58").Value = 2433.0527
$ws.Range("M58").Value = -563.4838999999999
$ws.Range("N58").Value = -2839.0527

# Row 100
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0

# Row 132
$ws.Range("H132").Value = 2200.5
$ws.Range("I132").Value = 1571.2653
$ws.Range("K132").Value = 4713.7959
$ws.Range("M132").Value = -2183.7959

# Row 136
$ws.Range("H136").Value = 1399.78
$ws.Range("I136").Value = 766.4838999999999
$ws.Range("J136").Value = 2433.0527
$ws.Range("K136").Value = 2299.4517
$ws.Range("L136").Value = 7299.158100000001
$ws.Range("M136").Value = 250.5483000000004
$ws.Range("N136").Value = -12399.1581

# Row 138
$ws.Range("H138").Value = 41196
$ws.Range("J138").Value = 41196
$ws.Range("L138").Value = 41196
$ws.Range("N138").Value = -51476

# Cleared cells (removed from source row)
$ws.Range("N100").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 17
$ws.Range("H17").Value = 771
$ws.Range("I17").Value = 41
$ws.Range("J17").Value = 1501
$ws.Range("K17").Value = 123
$ws.Range("L17").Value = 4503
$ws.Range("M17").Value = 46
$ws.Range("N17").Value = -4841

# Row 64
$ws.Range("H64").Value = 2787.4285
$ws.Range("J64").Value = 3140
$ws.Range("L64").Value = 9420
$ws.Range("N64").Value = -9960

# Row 67
$ws.Range("H67").Value = 2787.4285
$ws.Range("J67").Value = 3140
$ws.Range("L67").Value = 9420
$ws.Range("N67").Value = -11292

# Row 112
$ws.Range("H112").Value = 2647.0588
$ws.Range("I112").Value = 1750
$ws.Range("J112").Value = 2703.125
$ws.Range("K112").Value = 5250
$ws.Range("L112").Value = 8109.375
$ws.Range("M112").Value = -4142
$ws.Range("N112").Value = -10325.375

# Row 113
$ws.Range("H113").Value = 2000605.1
$ws.Range("I113").Value = 4546125.5
$ws.Range("J113").Value = 526882.8
$ws.Range("K113").Value = 13638376.5
$ws.Range("L113").Value = 1580648.4
$ws.Range("M113").Value = -13636206.5
$ws.Range("N113").Value = -1584988.4

# Row 119
$ws.Range("H119").Value = 145795.33
$ws.Range("I119").Value = 809.3333
$ws.Range("J119").Value = 290781.34
$ws.Range("K119").Value = 2427.9999
$ws.Range("L119").Value = 872344.02
$ws.Range("M119").Value = 2410.0001
$ws.Range("N119").Value = -882020.02

# Row 120
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0

# Row 121
$ws.Range("H121").Value = 993.902
$ws.Range("I121").Value = 938.46155
$ws.Range("J121").Value = 1012.8684
$ws.Range("K121").Value = 2815.38465
$ws.Range("L121").Value = 3038.6052
$ws.Range("M121").Value = -1505.38465
$ws.Range("N121").Value = -5658.6052

# Row 122
$ws.Range("H122").Value = 2578.7173
$ws.Range("I122").Value = 435.2
$ws.Range("J122").Value = 3615.9033
$ws.Range("K122").Value = 3916.8
$ws.Range("L122").Value = 32543.1297
$ws.Range("M122").Value = -1466.8
$ws.Range("N122").Value = -37443.1297

# Cleared cells (removed from source row)
$ws.Range("M120").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 122
$ws.Range("H122").Value = 22752628
$ws.Range("I122").Value = 38030270
$ws.Range("J122").Value = 2382435
$ws.Range("K122").Value = 114090810
$ws.Range("L122").Value = 7147305
$ws.Range("M122").Value = -114088360
$ws.Range("N122").Value = -7152205

# Row 132
$ws.Range("H132").Value = 2672.3438
$ws.Range("I132").Value = 2081.1052
$ws.Range("J132").Value = 3536.4614
$ws.Range("K132").Value = 6243.3156
$ws.Range("L132").Value = 10609.3842
$ws.Range("M132").Value = -3713.3156
$ws.Range("N132").Value = -15669.3842

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 132
$ws.Range("H132").Value = 13036764
$ws.Range("I132").Value = 14846734
$ws.Range("J132").Value = 4979.8
$ws.Range("K132").Value = 44540202
$ws.Range("L132").Value = 14939.4
$ws.Range("M132").Value = -44537672
$ws.Range("N132").Value = -19999.4

# Row 136
$ws.Range("H136").Value = 6953.1
$ws.Range("I136").Value = 4496.4
$ws.Range("J136").Value = 16779.9
$ws.Range("K136").Value = 13489.2
$ws.Range("L136").Value = 50339.7
$ws.Range("M136").Value = -10939.2
$ws.Range("N136").Value = -55439.7

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 41
$ws.Range("H41").Value = 8979.799999999999
$ws.Range("J41").Value = 8979.799999999999
$ws.Range("L41").Value = 8979.799999999999
$ws.Range("N41").Value = -9759.799999999999

# Row 113
$ws.Range("H113").Value = 787.71155
$ws.Range("I113").Value = 687.875
$ws.Range("J113").Value = 947.45
$ws.Range("K113").Value = 2063.625
$ws.Range("L113").Value = 2842.35
$ws.Range("M113").Value = 106.375
$ws.Range("N113").Value = -7182.35

# Row 132
$ws.Range("H132").Value = 17349.475
$ws.Range("I132").Value = 19220.648
$ws.Range("J132").Value = 2914.7144
$ws.Range("K132").Value = 57661.944
$ws.Range("L132").Value = 8744.143199999999
$ws.Range("M132").Value = -55131.944
$ws.Range("N132").Value = -13804.1432

# Row 136
$ws.Range("H136").Value = 10419760
$ws.Range("I136").Value = 4356.4585
$ws.Range("J136").Value = 20835164
$ws.Range("K136").Value = 13069.3755
$ws.Range("L136").Value = 62505492
$ws.Range("M136").Value = -10519.3755
$ws.Range("N136").Value = -62510592
